$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from paragraph 1
# ---------------------------------------------------------------------
try {
    $d.Bookmarks("_GoBack").Delete()
} catch {}

# ---------------------------------------------------------------------
# 2. Split paragraphs that need to be split into two, inserting a
#    leading space (copy-formatted from the preceding label's trailing
#    space) on the new paragraph.
# ---------------------------------------------------------------------

# VISTO: | La Ordenanza...
$d.Content.Find.Execute("VISTO: La Ordenanza", $true, $false, $false, $false, $false, $true, 1, $false, "VISTO: `rLa Ordenanza", 2) | Out-Null
$pVisto = $d.Paragraphs(3)
$pLa = $d.Paragraphs(4)
$srcSpace1 = $d.Range($pVisto.Range.End - 2, $pVisto.Range.End - 1)
$pLa.Range.InsertBefore(" ")
$newSpace1 = $d.Range($pLa.Range.Start, $pLa.Range.Start + 1)
$newSpace1.FormattedText = $srcSpace1.FormattedText

# CONSIDERANDO: | Que la mencionada...
$d.Content.Find.Execute("CONSIDERANDO: Que la mencionada", $true, $false, $false, $false, $false, $true, 1, $false, "CONSIDERANDO: `rQue la mencionada", 2) | Out-Null
$pConsid = $d.Paragraphs(5)
$pQue = $d.Paragraphs(6)
$srcSpace2 = $d.Range($pConsid.Range.End - 2, $pConsid.Range.End - 1)
$pQue.Range.InsertBefore(" ")
$newSpace2 = $d.Range($pQue.Range.Start, $pQue.Range.Start + 1)
$newSpace2.FormattedText = $srcSpace2.FormattedText

# ---------------------------------------------------------------------
# 3. Remove the leading "POR " run before "EL CONCEJO DELIBERANTE..."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("POR EL CONCEJO DELIBERANTE", $true, $false, $false, $false, $false, $true, 1, $false, "EL CONCEJO DELIBERANTE", 2) | Out-Null

# ---------------------------------------------------------------------
# After the two splits above the document now has 10 paragraphs:
#   1 Yerba Buena...
#   2 ORDENANZA Nº 1468
#   3 VISTO:
#   4  La Ordenanza...
#   5 CONSIDERANDO:
#   6  Que la mencionada...
#   7 EL CONCEJO DELIBERANTE...
#   8 ARTICULO PRIMERO: MODIFICASE...
#   9 "Nivel Sonoro...
#   10 ARTICULO SEGUNDO: COMUNIQUESE...
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 4. Paragraph-level direct formatting
# ---------------------------------------------------------------------

# P1 - Yerba Buena...
$p1 = $d.Paragraphs(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# P2 - ORDENANZA Nº 1468
$p2 = $d.Paragraphs(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = 1

# P3 - VISTO:
$p3 = $d.Paragraphs(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 12
$p3.Format.SpaceAfter = 6
$p3.Range.Font.Bold = 1

# P4 -  La Ordenanza...
$p4 = $d.Paragraphs(4)
$p4.Format.KeepWithNext = $true
$p4.Format.SpaceBefore = 0
$p4.Format.SpaceAfter = 6
$p4.Format.Alignment = 0

# P5 - CONSIDERANDO:
$p5 = $d.Paragraphs(5)
$p5.Format.KeepWithNext = $true
$p5.Format.SpaceBefore = 12
$p5.Format.SpaceAfter = 6
$p5.Format.Alignment = 0
$p5.Range.Font.Bold = 1

# P6 -  Que la mencionada...
$p6 = $d.Paragraphs(6)
$p6.Format.KeepWithNext = $true
$p6.Format.SpaceAfter = 6
$p6.Format.Alignment = 0

# P7 - EL CONCEJO DELIBERANTE...
$p7 = $d.Paragraphs(7)
$p7.Format.KeepWithNext = $true
$p7.Format.SpaceBefore = 18
$p7.Format.SpaceAfter = 18
$p7.Format.LeftIndent = 99.2
$p7.Format.RightIndent = 99.2
$p7.Range.Font.Bold = 1

# P8 - ARTICULO PRIMERO: MODIFICASE...
$p8 = $d.Paragraphs(8)
$p8.Format.KeepWithNext = $true
$p8.Format.SpaceAfter = 6
$p8.Format.Alignment = 0

# P9 - "Nivel Sonoro...
$p9 = $d.Paragraphs(9)
$p9.Format.KeepWithNext = $true
$p9.Format.SpaceAfter = 6
$p9.Format.Alignment = 0

# P10 - ARTICULO SEGUNDO: COMUNIQUESE...
$p10 = $d.Paragraphs(10)
$p10.Format.KeepWithNext = $true
$p10.Format.SpaceAfter = 6
$p10.Format.Alignment = 0

# ---------------------------------------------------------------------
# 5. Run-level direct formatting: underline "ARTICULO PRIMERO" / "SEGUNDO"
#    and split ": " into ":" (underlined) + " " (not underlined)
# ---------------------------------------------------------------------

# ARTICULO PRIMERO (paragraph 8)
$rngArt1 = $d.Range($p8.Range.Start, $p8.Range.End)
$rngArt1.Find.Execute("ARTICULO PRIMERO") | Out-Null
$rngArt1.Font.Underline = 1

$rngColon1 = $d.Range($p8.Range.Start, $p8.Range.End)
$rngColon1.Find.Execute(": ") | Out-Null
$colon1 = $d.Range($rngColon1.Start, $rngColon1.Start + 1)
$colon1.Font.Underline = 1

# ARTICULO SEGUNDO (paragraph 10)
$rngArt2 = $d.Range($p10.Range.Start, $p10.Range.End)
$rngArt2.Find.Execute("ARTICULO SEGUNDO") | Out-Null
$rngArt2.Font.Underline = 1

$rngColon2 = $d.Range($p10.Range.Start, $p10.Range.End)
$rngColon2.Find.Execute(": ") | Out-Null
$colon2 = $d.Range($rngColon2.Start, $rngColon2.Start + 1)
$colon2.Font.Underline = 1

# ---------------------------------------------------------------------
# 6. Section: footer + page numbering start
# ---------------------------------------------------------------------
$s = $d.Sections(1)
$f = $s.Footers(1)
$fp = $f.Range.Paragraphs(1)
$fp.Style = "Piedepgina"
$fp.Range.Font.Name = "Book Antiqua"
$fp.Range.Font.Size = 10
$fp.Range.Font.Color = 8421504

$f.PageNumbers.StartingNumber = 1681

# header / footer paragraph styles referenced by the style gallery even
# though no header content is added
$stH = $d.Styles.Add("Encabezado", 1)
$stHC = $d.Styles.Add("EncabezadoCar", 2)
$stHC.NameLocal = "Encabezado Car"
$stH.BaseStyle = "Normal"
$stH.LinkStyle = "EncabezadoCar"
$stHC.LinkStyle = "Encabezado"
$stH.Priority = 99
$stHC.Priority = 99
$stH.UnhideWhenUsed = $true

$stF = $d.Styles("Piedepgina")
$stFC = $d.Styles.Add("PiedepginaCar", 2)
$stFC.NameLocal = "Pie de página Car"
$stF.LinkStyle = "PiedepginaCar"
$stFC.LinkStyle = "Piedepgina"
$stF.Priority = 99
$stFC.Priority = 99
$stF.UnhideWhenUsed = $true

# ---------------------------------------------------------------------
# 7. Footnote / endnote separator parts (touched so the default
#    separator/continuationSeparator markup is minted)
# ---------------------------------------------------------------------
$fn = $d.Footnotes.Add($d.Paragraphs(1).Range, "", "x")
$fn.Delete()

Write-Output "done"
